# Update "想去人数" (want-to-go count, column F) figures across the four
# sheets to match the refreshed data snapshot ("output generated at 456a3b4").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 453
$ws.Range("F8").Value = 1189
$ws.Range("F9").Value = 339
$ws.Range("F11").Value = 877
$ws.Range("F12").Value = 682
$ws.Range("F14").Value = 505
$ws.Range("F17").Value = 173
$ws.Range("F18").Value = 2919
$ws.Range("F20").Value = 525
$ws.Range("F24").Value = 227
$ws.Range("F26").Value = 5268
$ws.Range("F27").Value = 590
$ws.Range("F28").Value = 983
$ws.Range("F29").Value = 22
$ws.Range("F31").Value = 308
$ws.Range("F32").Value = 1094
$ws.Range("F34").Value = 54
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1118
$ws.Range("F5").Value = 6
$ws.Range("F10").Value = 29
$ws.Range("F17").Value = 986
$ws.Range("F26").Value = 3924
$ws.Range("F29").Value = 21
$ws.Range("F33").Value = 161
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2453
$ws.Range("F6").Value = 1040
$ws.Range("F9").Value = 1320
$ws.Range("F10").Value = 359
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2453
$ws.Range("F6").Value = 1040
$ws.Range("F7").Value = 1320
$ws.Range("F8").Value = 359
$ws.Range("F11").Value = 453
$ws.Range("F14").Value = 1189
$ws.Range("F15").Value = 339
$ws.Range("F16").Value = 877
$ws.Range("F17").Value = 682
$ws.Range("F18").Value = 1118
$ws.Range("F19").Value = 1118
$ws.Range("F20").Value = 505
$ws.Range("F22").Value = 173
$ws.Range("F23").Value = 2919
$ws.Range("F25").Value = 525
$ws.Range("F28").Value = 29
$ws.Range("F29").Value = 227
$ws.Range("F30").Value = 5268
$ws.Range("F31").Value = 590
$ws.Range("F32").Value = 983
$ws.Range("F35").Value = 22
$ws.Range("F38").Value = 308
$ws.Range("F47").Value = 1094
$ws.Range("F49").Value = 161
$ws.Range("F50").Value = 54
